$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.957.86'
$ws.Range('E2').Value = '  -1.99%  '
$ws.Range('D3').Value = '3.106.67'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '526.14'
$ws.Range('E5').Value = '  +0.34%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.19'
$ws.Range('E6').Value = '  -2.10%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '3.106.01'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.444'
$ws.Range('E9').Value = '  +1.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.16'
$ws.Range('E10').Value = '  -3.27%  '
$ws.Range('E11').Value = '  -1.46%  '
$ws.Range('E12').Value = '  +1.94%  '
$ws.Range('D13').Value = '3.641.38'
$ws.Range('E13').Value = '  -0.05%  '
$ws.Range('E14').Value = '  +3.20%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.63'
$ws.Range('E15').Value = '  -5.46%  '
$ws.Range('E16').Value = '  -1.40%  '
$ws.Range('D17').Value = '58.000.89'
$ws.Range('E17').Value = '  -1.91%  '
$ws.Range('D18').Value = '3.103.66'
$ws.Range('E18').Value = '  -0.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.12'
$ws.Range('E19').Value = '  -1.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.72'
$ws.Range('E20').Value = '  -2.52%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.97'
$ws.Range('E21').Value = '  -2.71%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '342.88'
$ws.Range('E22').Value = '  -0.34%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.513'
$ws.Range('E24').Value = '  +0.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '67.55'
$ws.Range('E25').Value = '  +2.45%  '
$ws.Range('E26').Value = '  -0.66%  '
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('D28').Value = '0.0₃0928'
$ws.Range('E28').Value = '  -0.93%  '
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.29'
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('B31').Value = 'RenderToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.35'
$ws.Range('E31').Value = '  -6.76%  '
$ws.Range('E32').Value = '  +1.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.03'
$ws.Range('E33').Value = '  -0.17%  '
$ws.Range('E34').Value = '  -3.11%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '158.60'
$ws.Range('E35').Value = '  +2.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.63'
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('E37').Value = '  -0.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '26.24'
$ws.Range('E38').Value = '  -2.96%  '
$ws.Range('E39').Value = '  -4.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0669'
$ws.Range('E40').Value = '  -2.77%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.62'
$ws.Range('E41').Value = '  +10.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.00'
$ws.Range('E42').Value = '  +0.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.688'
$ws.Range('E43').Value = '  +3.15%  '
$ws.Range('D44').Value = '3.145.69'
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '36.82'
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0262'
$ws.Range('E47').Value = '  +1.46%  '
$ws.Range('D48').Value = '2.270.35'
$ws.Range('E48').Value = '  -1.50%  '
$ws.Range('E49').Value = '  +2.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.14'
$ws.Range('E50').Value = '  +1.55%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.54'
$ws.Range('E51').Value = '  -1.84%  '
